# cluster_table.xlsx - data sanitizing update
# Re-label several "Functional Type" entries (column E of the cluster_info
# sheet / Table5) from their prior category to the newly-introduced "Mixed"
# category, and scroll the sheet view back up a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cluster_info")

# Rows whose Functional Type (column E) is reclassified to "Mixed" as part
# of the data-sanitizing pass described in the commit message.
$rowsToMix = @(22, 23, 25, 44, 51, 53, 58, 62, 65, 67, 70, 71, 72, 83, 84, 85, 88, 90, 91)

foreach ($r in $rowsToMix) {
    $ws.Range("E$r").Value = "Mixed"
}

# Scroll the sheet so row 69 is the top visible row (was row 85), and
# restore the original selection.
$ws.Range("A69").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 69
$excel.ActiveWindow.ScrollColumn = 1

$ws.Range("A1:H96").Select() | Out-Null
